$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

# Copy the date-cell formatting from the row above so the new date cell
# reuses the existing style (no new style entries get created).
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 42611.885659722226
$ws.Cells.Item($row, 2).Value = 2
$ws.Cells.Item($row, 3).Value = 52
$ws.Cells.Item($row, 4).Value = 43
$ws.Cells.Item($row, 5).Value = 47
$ws.Cells.Item($row, 6).Value = 52
$ws.Cells.Item($row, 7).Value = 10427
$ws.Cells.Item($row, 8).Value = 19550
$ws.Cells.Item($row, 9).Value = 2478
$ws.Cells.Item($row, 10).Value = 239
$ws.Cells.Item($row, 11).Value = 198
$ws.Cells.Item($row, 12).Value = 11
$ws.Cells.Item($row, 13).Value = 12
$ws.Cells.Item($row, 14).Value = "Bag"
